$wb = $excel.ActiveWorkbook

$wsKim = $wb.Worksheets.Item("Kim")
$wsSam = $wb.Worksheets.Item("Sam")

# --- Sheet "Kim": insert 5 new year rows (2050-2054) right before the
#     trailing 2055/2056/2057 rows, then renumber the whole year column
#     so it now starts 5 years earlier (2020 instead of 2025). ---
$wsKim.Rows("32:36").Insert()
for ($i = 0; $i -lt 38; $i++) {
    $wsKim.Cells.Item(2 + $i, 1).Value = 2020 + $i
}

# --- Sheet "Sam": same change, inserted right before its trailing
#     2057-only row. ---
$wsSam.Rows("34:38").Insert()
for ($i = 0; $i -lt 38; $i++) {
    $wsSam.Cells.Item(2 + $i, 1).Value = 2020 + $i
}

# --- Selections / active sheet, matching the saved view state ---
$wsKim.Range("A2:XFD6").Select() | Out-Null
$wsSam.Activate() | Out-Null
$wsSam.Range("B9").Select() | Out-Null
